$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D (price) to Text format so that purely numeric-looking
# price strings (e.g. "187.24") are written back as text, matching the original
# inlineStr cell type, instead of being auto-converted to floating point numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "65.799.54"
$ws.Range("E2").Value = "  +1.01%  "

$ws.Range("D3").Value = "3.310.22"
$ws.Range("E3").Value = "  +1.38%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "187.24"
$ws.Range("E5").Value = "  +5.23%  "

$ws.Range("D6").Value = "553.93"
$ws.Range("E6").Value = "  +0.44%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "3.305.70"
$ws.Range("E8").Value = "  +1.34%  "

$ws.Range("D9").Value = "0.579"
$ws.Range("E9").Value = "  -0.80%  "

$ws.Range("D10").Value = "0.179"
$ws.Range("E10").Value = "  -3.14%  "

$ws.Range("D11").Value = "0.580"
$ws.Range("E11").Value = "  -0.05%  "

$ws.Range("D12").Value = "46.68"
$ws.Range("E12").Value = "  -0.41%  "

$ws.Range("D13").Value = "0.0000266"
$ws.Range("E13").Value = "  +2.12%  "

$ws.Range("D14").Value = "8.61"
$ws.Range("E14").Value = "  +1.68%  "

$ws.Range("D15").Value = "3.833.16"
$ws.Range("E15").Value = "  +1.42%  "

$ws.Range("D16").Value = "597.36"
$ws.Range("E16").Value = "  -0.89%  "

$ws.Range("D17").Value = "65.815.16"
$ws.Range("E17").Value = "  +1.11%  "

$ws.Range("D19").Value = "17.86"
$ws.Range("E19").Value = "  +0.45%  "

$ws.Range("D20").Value = "3.306.68"
$ws.Range("E20").Value = "  +1.70%  "

$ws.Range("D21").Value = "10.98"
$ws.Range("E21").Value = "  -2.94%  "

$ws.Range("D22").Value = "0.896"
$ws.Range("E22").Value = "  +0.28%  "

$ws.Range("D23").Value = "18.46"
$ws.Range("E23").Value = "  +6.08%  "

$ws.Range("D24").Value = "5.04"
$ws.Range("E24").Value = "  +1.83%  "

$ws.Range("D25").Value = "100.21"
$ws.Range("E25").Value = "  -1.17%  "

$ws.Range("E26").Value = "  -0.32%  "

$ws.Range("E27").Value = "  +1.29%  "

$ws.Range("D28").Value = "2.73"
$ws.Range("E28").Value = "  +2.85%  "

$ws.Range("D29").Value = "9.44"
$ws.Range("E29").Value = "  +1.71%  "

$ws.Range("D30").Value = "8.66"
$ws.Range("E30").Value = "  +0.98%  "

$ws.Range("D31").Value = "30.33"
$ws.Range("E31").Value = "  +0.50%  "

$ws.Range("D32").Value = "6.67"
$ws.Range("E32").Value = "  +8.00%  "

$ws.Range("D33").Value = "3.84"
$ws.Range("E33").Value = "  +0.78%  "

$ws.Range("D34").Value = "568.47"
$ws.Range("E34").Value = "  +5.98%  "

$ws.Range("D35").Value = "10.98"
$ws.Range("E35").Value = "  +0.32%  "

$ws.Range("E36").Value = "  -0.14%  "

$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "3.705.71"
$ws.Range("E37").Value = "  -0.68%  "

$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.03%  "

$ws.Range("D39").Value = "56.67"
$ws.Range("E39").Value = "  +1.26%  "

$ws.Range("D40").Value = "3.52"
$ws.Range("E40").Value = "  +11.05%  "

$ws.Range("D41").Value = "33.52"
$ws.Range("E41").Value = "  +5.79%  "

$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "0.127"
$ws.Range("E42").Value = "  +2.23%  "

$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "3.21"
$ws.Range("E43").Value = "  -5.98%  "

$ws.Range("D44").Value = "0.0₃0696"
$ws.Range("E44").Value = "  -0.31%  "

$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").Value = "2.64"
$ws.Range("E45").Value = "  -0.46%  "

$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "3.42"
$ws.Range("E46").Value = "  +8.31%  "

$ws.Range("E47").Value = "  +0.11%  "

$ws.Range("D48").Value = "0.0414"
$ws.Range("E48").Value = "  +2.11%  "

$ws.Range("D49").Value = "0.129"
$ws.Range("E49").Value = "  +0.29%  "

$ws.Range("E50").Value = "  +0.15%  "

$ws.Range("E51").Value = "  -0.73%  "

# Restore the default (Normal) style on column D so no stray number-format
# styling is left behind on cells that did not originally have one.
$priceRange.Style = "Normal"
